# Scheduled data refresh: update market-board price/profit columns (H:N)
# on the Leve profit-tracking sheets. Values come from the latest
# currentAveragePrice* pull; LeveProfit* (M/N) are recomputed accordingly.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 140.75
$ws.Range("I5").Value = 123.4
$ws.Range("J5").Value = 169.66667
$ws.Range("K5").Value = 123.4
$ws.Range("L5").Value = 169.66667
$ws.Range("M5").Value = -8.400000000000006
$ws.Range("N5").Value = -399.66667

$ws.Range("H43").Value = 2859.8333
$ws.Range("I43").Value = 1559
$ws.Range("J43").Value = 3120
$ws.Range("K43").Value = 1559
$ws.Range("L43").Value = 3120
$ws.Range("M43").Value = -1490
$ws.Range("N43").Value = -3258

$ws.Range("H74").Value = 22093.889
$ws.Range("I74").Value = 24931.5
$ws.Range("J74").Value = 16418.666
$ws.Range("K74").Value = 24931.5
$ws.Range("L74").Value = 16418.666
$ws.Range("M74").Value = -23995.5
$ws.Range("N74").Value = -18290.666

$ws.Range("H77").Value = 22093.889
$ws.Range("I77").Value = 24931.5
$ws.Range("J77").Value = 16418.666
$ws.Range("K77").Value = 124657.5
$ws.Range("L77").Value = 82093.33
$ws.Range("M77").Value = -119977.5
$ws.Range("N77").Value = -91453.33

$ws.Range("H93").Value = 30000
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 30000
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 30000
$ws.Range("N93").Value = -34992

$ws.Range("H100").Value = 2998.7144
$ws.Range("I100").Value = 2598.4
$ws.Range("J100").Value = 3999.5
$ws.Range("K100").Value = 2598.4
$ws.Range("L100").Value = 3999.5
$ws.Range("M100").Value = -2057.4

$ws.Range("H132").Value = 1738.4
$ws.Range("I132").Value = 1755.5
$ws.Range("J132").Value = 1499
$ws.Range("K132").Value = 5266.5
$ws.Range("L132").Value = 4497
$ws.Range("M132").Value = -2736.5
$ws.Range("N132").Value = -9557

$ws.Range("H137").Value = 1198.1428
$ws.Range("I137").Value = 1189.25
$ws.Range("J137").Value = 1210
$ws.Range("K137").Value = 3567.75
$ws.Range("L137").Value = 3630
$ws.Range("M137").Value = -1017.75
$ws.Range("N137").Value = -8730

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2028659.4
$ws.Range("I32").Value = 2416228.5
$ws.Range("J32").Value = 779825.25
$ws.Range("K32").Value = 2416228.5
$ws.Range("L32").Value = 779825.25
$ws.Range("M32").Value = -2415941.5
$ws.Range("N32").Value = -780399.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 22573.666
$ws.Range("I7").Value = 1055
$ws.Range("J7").Value = 33333
$ws.Range("K7").Value = 1055
$ws.Range("L7").Value = 33333
$ws.Range("M7").Value = -942
$ws.Range("N7").Value = -33559

$ws.Range("H107").Value = 3737.125
$ws.Range("I107").Value = 3296.5
$ws.Range("J107").Value = 5059
$ws.Range("K107").Value = 3296.5
$ws.Range("L107").Value = 5059
$ws.Range("M107").Value = -1376.5
$ws.Range("N107").Value = -8899

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1200
$ws.Range("I31").Value = 1200
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1200
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -905

$ws.Range("H34").Value = 1200
$ws.Range("I34").Value = 1200
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1200
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -998

$ws.Range("H58").Value = 2032.7142
$ws.Range("I58").Value = 2002.4
$ws.Range("J58").Value = 2108.5
$ws.Range("K58").Value = 2002.4
$ws.Range("L58").Value = 2108.5
$ws.Range("M58").Value = -1799.4
$ws.Range("N58").Value = -2514.5

$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()

$ws.Range("H105").Value = 3200
$ws.Range("I105").Value = 2500
$ws.Range("J105").Value = 3410
$ws.Range("K105").Value = 2500
$ws.Range("L105").Value = 3410
$ws.Range("M105").Value = -753

$ws.Range("H132").Value = 3495
$ws.Range("I132").Value = 3495
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 10485
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -7955
$ws.Range("N132").ClearContents()

$ws.Range("H134").Value = 2113.875
$ws.Range("I134").Value = 2068.6667
$ws.Range("J134").Value = 2249.5
$ws.Range("K134").Value = 6206.000100000001
$ws.Range("L134").Value = 6748.5
$ws.Range("M134").Value = -3671.000100000001
$ws.Range("N134").Value = -11818.5

$ws.Range("H136").Value = 2032.7142
$ws.Range("I136").Value = 2002.4
$ws.Range("J136").Value = 2108.5
$ws.Range("K136").Value = 6007.200000000001
$ws.Range("L136").Value = 6325.5
$ws.Range("M136").Value = -3457.200000000001
$ws.Range("N136").Value = -11425.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 1000
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 3000
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -2831

$ws.Range("H27").Value = 1000
$ws.Range("I27").Value = 1000
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 3000
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -2898

$ws.Range("H37").Value = 69988.5
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 69988.5
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 209965.5
$ws.Range("N37").Value = -210189.5

$ws.Range("H86").Value = 968.6667
$ws.Range("I86").Value = 900
$ws.Range("J86").Value = 1003
$ws.Range("K86").Value = 2700
$ws.Range("L86").Value = 3009
$ws.Range("M86").Value = -1514
$ws.Range("N86").Value = -5381

$ws.Range("H89").Value = 968.6667
$ws.Range("I89").Value = 900
$ws.Range("J89").Value = 1003
$ws.Range("K89").Value = 8100
$ws.Range("L89").Value = 9027
$ws.Range("M89").Value = -2172
$ws.Range("N89").Value = -20883

$ws.Range("H113").Value = 1578.4445
$ws.Range("I113").Value = 1405.5
$ws.Range("J113").Value = 1716.8
$ws.Range("K113").Value = 4216.5
$ws.Range("L113").Value = 5150.4
$ws.Range("M113").Value = -2046.5
$ws.Range("N113").Value = -9490.4

$ws.Range("H121").Value = 3089.2778
$ws.Range("I121").Value = 2640
$ws.Range("J121").Value = 3145.4375
$ws.Range("K121").Value = 7920
$ws.Range("L121").Value = 9436.3125
$ws.Range("M121").Value = -6610
$ws.Range("N121").Value = -12056.3125

$ws.Range("H128").Value = 612088.25
$ws.Range("I128").Value = 612088.25
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 1836264.75
$ws.Range("L128").Value = 0
$ws.Range("M128").Value = -1831284.75

$ws.Range("H132").Value = 6349
$ws.Range("I132").Value = 4999.5
$ws.Range("J132").Value = 7698.5
$ws.Range("K132").Value = 44995.5
$ws.Range("L132").Value = 69286.5
$ws.Range("M132").Value = -42465.5
$ws.Range("N132").Value = -74346.5

$ws.Range("H134").Value = 1966.6666
$ws.Range("I134").Value = 1966.6666
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 5899.9998
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -829.9997999999996

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3716.625
$ws.Range("I80").Value = 3777.6667
$ws.Range("J80").Value = 3680
$ws.Range("K80").Value = 3777.6667
$ws.Range("L80").Value = 3680
$ws.Range("M80").Value = -2779.6667

$ws.Range("H83").Value = 3716.625
$ws.Range("I83").Value = 3777.6667
$ws.Range("J83").Value = 3680
$ws.Range("K83").Value = 18888.3335
$ws.Range("L83").Value = 18400
$ws.Range("M83").Value = -13896.3335

$ws.Range("H134").Value = 41272
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 41272
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 123816
$ws.Range("N134").Value = -128886

$ws.Range("H136").Value = 48774.668
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 48774.668
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 146324.004
$ws.Range("N136").Value = -151424.004

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2148.7273
$ws.Range("I136").Value = 1515.2222
$ws.Range("J136").Value = 4999.5
$ws.Range("K136").Value = 4545.6666
$ws.Range("L136").Value = 14998.5
$ws.Range("M136").Value = -1995.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 19976
$ws.Range("I41").Value = 19978
$ws.Range("J41").Value = 19975.5
$ws.Range("K41").Value = 19978
$ws.Range("L41").Value = 19975.5
$ws.Range("M41").Value = -19588
$ws.Range("N41").Value = -20755.5

$ws.Range("H132").Value = 7666.6665
$ws.Range("I132").Value = 7666.6665
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 22999.9995
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -20469.9995
$ws.Range("N132").ClearContents()

$ws.Range("H136").Value = 2344.842
$ws.Range("I136").Value = 2239.7144
$ws.Range("J136").Value = 2639.2
$ws.Range("K136").Value = 6719.1432
$ws.Range("L136").Value = 7917.599999999999
$ws.Range("M136").Value = -4169.1432
$ws.Range("N136").Value = -13017.6
